$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new translator name "Experiment2" across T1:AE1 ---
$ws.Range("B1:M1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("T1:AE1").Value = "published_Colloff_Flowe_SealeCarlisle_2020_Experiment2()"

# --- Row 2: model parameter combo, repeats the B2:G2 pattern twice ---
$ws.Range("B2:M2").Copy()
$ws.Range("T2").PasteSpecial(-4122)
$ws.Range("B2:G2").Copy()
$ws.Range("T2").PasteSpecial(-4163)
$ws.Range("B2:G2").Copy()
$ws.Range("Z2").PasteSpecial(-4163)

# --- Row 3: model parameter combo, repeats the B3:G3 pattern twice ---
$ws.Range("B3:M3").Copy()
$ws.Range("T3").PasteSpecial(-4122)
$ws.Range("B3:G3").Copy()
$ws.Range("T3").PasteSpecial(-4163)
$ws.Range("B3:G3").Copy()
$ws.Range("Z3").PasteSpecial(-4163)

# --- Row 4: exclusions "{}" throughout ---
$ws.Range("B4:M4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4:AE4").Value = "{}"

# --- Row 5: new encoding-condition labels ---
$ws.Range("B5:M5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5:Y5").Value = "encodingConditionRaw 1"
$ws.Range("Z5:AE5").Value = "encodingConditionRaw 2"

# --- Row 6: binning "[-1,60,80,100]" throughout ---
$ws.Range("B6:M6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6:AE6").Value = "[-1,60,80,100]"

# --- Row 7: "True" throughout (copy the whole cell so the text stays a string, not boolean) ---
$ws.Range("B7:M7").Copy()
$ws.Range("T7").PasteSpecial()

# --- Row 8: niter 2000 throughout ---
$ws.Range("B8:M8").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$ws.Range("T8:AE8").Value = 2000

$excel.CutCopyMode = 0

# --- match the "bestFit" 35-wide columns used by B:S, now extended through AE ---
$ws.Columns("T:AE").ColumnWidth = $ws.Range("B1").ColumnWidth

# --- sheet view / selection housekeeping to mirror the authored edit ---
$ws.Range("AA5:AE5").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 24
